$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45565
$ws.Range("B2").Value = "SEI-260007/016010/2024"
$ws.Range("C2").Value = "DIVLS"
$ws.Range("D2").Value = "Licitação"
$ws.Range("E2").ClearContents()

$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "1"

$ws.Range("G2").Value = "Maria Eduarda"
$ws.Range("H2").Value = 45565
$ws.Range("I2").Value = "10034/2024"
$ws.Range("J2").Value = 45666
$ws.Range("M2").Value = 45667
$ws.Range("N2").Value = 102
$ws.Range("O2").Value = "September"

$ws.Range("P2").NumberFormat = "@"
$ws.Range("P2").Value = "1"

$ws.Range("Q2").Value = 365
$ws.Range("R2").Value = "Saiu em 10/01/2025"
